$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows above the current row 5 (old rows 5 "Chad" and 6 "Aland"
# shift down to become rows 11 and 12). Looping single-row inserts keeps the
# row-by-row formatting inheritance (s=2 / s=3 on column D) the same way
# Excel does it when a user inserts rows one at a time.
for ($i = 0; $i -lt 6; $i++) {
    $ws.Rows.Item(5).Insert()
}

# Row 4 gained a source_id value that used to be blank.
$ws.Cells.Item(4, 7).Value = 1

# The three new "preview citation" sentences become new shared-string table
# entries the first time each distinct text is written. Write them in
# Lincoln / Nebraska / United States order (matching the first appearance
# order a full regenerate-preview pass would hit) so new shared-string
# indices land the same way as the target file.
$ws.Cells.Item(6, 9).Value = "philbert in Lincoln as described by Anon."
$ws.Cells.Item(5, 9).Value = "philbert in Nebraska as described by Anon."
$ws.Cells.Item(8, 9).Value = "philbert in the United States as described by Anon."
$ws.Cells.Item(7, 9).Value = "philbert in Lincoln as described by Anon."
$ws.Cells.Item(9, 9).Value = "philbert in Nebraska as described by Anon."
$ws.Cells.Item(10, 9).Value = "philbert in Lincoln as described by Anon."

# Row 5 - philbert / United States / Nebraska / source 2
$ws.Cells.Item(5, 1).Value = "philbert"
$ws.Cells.Item(5, 4).Value = "United States"
$ws.Cells.Item(5, 5).Value = "Nebraska"
$ws.Cells.Item(5, 7).Value = 2
$ws.Rows.Item(5).RowHeight = 75

# Row 6 - philbert / United States / Lincoln / source 3
$ws.Cells.Item(6, 1).Value = "philbert"
$ws.Cells.Item(6, 4).Value = "United States"
$ws.Cells.Item(6, 6).Value = "Lincoln"
$ws.Cells.Item(6, 7).Value = 3
$ws.Rows.Item(6).RowHeight = 60

# Row 7 - philbert / Nebraska / Lincoln / source 4 (no country this time)
$ws.Cells.Item(7, 1).Value = "philbert"
$ws.Cells.Item(7, 5).Value = "Nebraska"
$ws.Cells.Item(7, 6).Value = "Lincoln"
$ws.Cells.Item(7, 7).Value = 4
$ws.Rows.Item(7).RowHeight = 60

# Row 8 - philbert / United States / source 5
$ws.Cells.Item(8, 1).Value = "philbert"
$ws.Cells.Item(8, 4).Value = "United States"
$ws.Cells.Item(8, 7).Value = 5
$ws.Rows.Item(8).RowHeight = 75

# Row 9 - philbert / Nebraska / source 9
$ws.Cells.Item(9, 1).Value = "philbert"
$ws.Cells.Item(9, 5).Value = "Nebraska"
$ws.Cells.Item(9, 7).Value = 9
$ws.Rows.Item(9).RowHeight = 75

# Row 10 - philbert / Lincoln / source 10
$ws.Cells.Item(10, 1).Value = "philbert"
$ws.Cells.Item(10, 6).Value = "Lincoln"
$ws.Cells.Item(10, 7).Value = 10
$ws.Rows.Item(10).RowHeight = 60

# Give the preview-text cells (I5:I9) the look of pasted-in generated text:
# explicit black font colour + vertical-centred wrap, matching the new
# cellXfs style introduced for this batch-preview column. Row 10's I cell
# keeps the plain default style.
$previewRange = $ws.Range("I5:I9")
$previewRange.Font.Color = 0
$previewRange.WrapText = $true
$previewRange.VerticalAlignment = -4108

# Reselect / rescroll to show the newly added preview rows, like a user
# would after generating the preview.
$ws.Range("I6").Select()

# Sheet now prints in portrait orientation.
$ws.PageSetup.Orientation = 1
